# Update column G ("K" - strikeouts) values for rows 2-34 on Sheet1.
# This regenerates the save_data to use actual strikeout counts (K)
# instead of the previous "Strike#" (total pitch/strike count) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 5
    3  = 2
    4  = 2
    5  = 4
    6  = 2
    7  = 3
    8  = 6
    9  = 8
    10 = 2
    11 = 0
    12 = 9
    13 = 3
    14 = 5
    15 = 4
    16 = 4
    17 = 5
    18 = 3
    19 = 4
    20 = 4
    21 = 6
    22 = 3
    23 = 7
    24 = 1
    25 = 5
    26 = 0
    27 = 2
    28 = 4
    29 = 5
    30 = 2
    31 = 2
    32 = 6
    33 = 3
    34 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
